$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Output $ws.Name
